$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrls = @(
    'https://cryptojobslist.com/jobs/head-of-engineering-cere-network-europe',
    'https://cryptojobslist.com/jobs/senior-backend-engineer-dappradar-remote-1',
    'https://cryptojobslist.com/jobs/incident-responder-apac-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/principal-engineer-chainlink-labs-remote',
    'https://cryptojobslist.com/marketing/market-strategy-manager-capital-markets-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/senior-data-partnerships-manager-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/strategic-finance-manager-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/research-engineer-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/senior-site-reliability-engineer-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/technical-account-manager-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/engineering-manager-blockchain-integrations-chainlink-labs-remote',
    'https://cryptojobslist.com/jobs/senior-people-business-partner-chainlink-labs-remote',
    'https://cryptojobslist.com/marketing/social-media-marketing-specialist-apex-digital-solutions-dubai',
    'https://cryptojobslist.com/marketing/brand-strategist-copywriter-otherlife-moonpay-toronto-canada',
    'https://cryptojobslist.com/jobs/director-of-total-rewards-moonpay-uk-eu',
    'https://cryptojobslist.com/jobs/internship-development-certik-new-york-seattle-bay-area-remote',
    'https://cryptojobslist.com/jobs/devops-engineer-art-blocks-remote',
    'https://cryptojobslist.com/jobs/head-of-infrastructure-chainsafe-systems-remote',
    'https://cryptojobslist.com/marketing/head-of-sales-chainsafe-systems-toronto',
    'https://cryptojobslist.com/jobs/senior-software-engineer-frontend-chainalysis-remote-can',
    'https://cryptojobslist.com/jobs/front-end-engineer-aptos-palo-alto',
    'https://cryptojobslist.com/jobs/region-lead-south-east-asia-scroll-io-singapore',
    'https://cryptojobslist.com/marketing/developer-relations-sygma-chainsafe-systems-toronto-canada-or-remote',
    'https://cryptojobslist.com/jobs/senior-agile-practitioner-chainsafe-systems-remote',
    'https://cryptojobslist.com/jobs/executive-assistant-to-the-ceo-subspace-labs-palo-alto-ca',
    'https://cryptojobslist.com/jobs/kyc-analyst-japan-binance-asia',
    'https://cryptojobslist.com/jobs/wallet-operations-specialist-binance-asia',
    'https://cryptojobslist.com/jobs/vp-communications-yuga-labs-remote',
    'https://cryptojobslist.com/jobs/compliance-officer-money-laundering-reporting-officer-co-mlro-measa-binance-uae-dubai',
    'https://cryptojobslist.com/jobs/bap-spot-wallet-operations-binance-asia',
    'https://cryptojobslist.com/jobs/suspicious-activity-report-sar-investigator-binance-asia',
    'https://cryptojobslist.com/jobs/director-of-core-compliance-programs-binance-ireland-dublin',
    'https://cryptojobslist.com/marketing/senior-product-manager-gaming-chainsafe-systems-berlin-or-remote-cet-est-hours-only',
    'https://cryptojobslist.com/jobs/senior-backend-engineer-uniswap-labs-new-york-or-remote',
    'https://cryptojobslist.com/marketing/growth-marketing-manager-uniswap-labs-new-york-new-york-united-states',
    'https://cryptojobslist.com/marketing/africa-senior-marketing-manager-binance-south-africa-cape-town',
    'https://cryptojobslist.com/marketing/graphic-designer-remote-metaschool-singapore',
    'https://cryptojobslist.com/jobs/global-head-of-business-development-amp-partnerships-remote-metaschool-singapore',
    'https://cryptojobslist.com/jobs/sales-manager-metaschool-singapore'
)

$startRow = 102
for ($i = 0; $i -lt $newUrls.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newUrls[$i]
}
